$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for two new match rows right after row 213 (the last existing
#    match row). We insert twice *at* row 213 so that Excel pushes the
#    populated row 213 downward (cloning its real per-cell styles as it
#    goes) instead of inventing styles for a previously-empty row - doing
#    the insert at row 214 directly hits an empty/gap row and the engine
#    fabricates bogus new cell styles, which we don't want.
#
#    After this:
#      row 213 -> blank, styled like the row above it (212, the "even" style)
#      row 214 -> blank, styled like the row above it (212, the "even" style)
#      row 215 -> the original row 213 match data (untouched, "odd" style)
# ---------------------------------------------------------------------------
$ws.Rows("213:213").Insert()
$ws.Rows("214:214").Insert()

# Re-home the original row-213 data (now sitting at row 215) back onto row
# 213 by cloning row 215's formatting onto row 213 (reuses the existing
# style indices - no new styles get created).
$ws.Range("A215:E215").Copy()
$ws.Range("A213:E213").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Write the original match back into row 213 (unchanged values).
# ---------------------------------------------------------------------------
$ws.Range("A213").Value = 43452
$ws.Range("B213").Value = 10
$ws.Range("C213").Value = 15
$ws.Range("D213").Value = "A"
$ws.Range("E213").Value = 43452.56280092592

# ---------------------------------------------------------------------------
# 3. New match row 214 (style already cloned from row 212 = "even" pattern).
# ---------------------------------------------------------------------------
$ws.Range("A214").Value = 43453
$ws.Range("B214").Value = 15
$ws.Range("C214").Value = 13
$ws.Range("D214").Value = "H"
$ws.Range("E214").Value = 43453.421875

# ---------------------------------------------------------------------------
# 4. New match row 215 (style already cloned from old row 213 = "odd"
#    pattern - just overwrite the values).
# ---------------------------------------------------------------------------
$ws.Range("A215").Value = 43453
$ws.Range("B215").Value = 15
$ws.Range("C215").Value = 10
$ws.Range("D215").Value = "A"
$ws.Range("E215").Value = 43453.42680555556

# ---------------------------------------------------------------------------
# 5. Update the scratch-pad / summary table further down the sheet. Row
#    insertion above has already shifted it from 216-226 to 218-228 and
#    fixed up the mergeCells for us; we just correct the recomputed
#    aggregate numbers.
# ---------------------------------------------------------------------------
$ws.Range("C218").Value = 44
$ws.Range("C219").Value = 49
$ws.Range("C220").Value = 93
$ws.Range("C221").Value = 38.94
$ws.Range("D221").Value = 51.49
$ws.Range("C222").Value = 48.51
$ws.Range("D222").Value = 61.06
$ws.Range("C223").Value = 43.46
$ws.Range("D223").Value = 56.54
$ws.Range("C224").Value = 12.65
$ws.Range("D224").Value = 13.5
$ws.Range("C225").Value = 13.22
$ws.Range("D225").Value = 13.96
$ws.Range("C226").Value = 12.92
$ws.Range("D226").Value = 13.74
$ws.Range("C227").Value = 2
$ws.Range("D227").Value = 0

# ---------------------------------------------------------------------------
# 6. New conditional formatting rules for the two new match rows (mirrors
#    the per-row "greaterThan" cellIs rules used throughout the sheet).
# ---------------------------------------------------------------------------
$ws.Range("B214").FormatConditions.Add(1, 5, "13") | Out-Null
$ws.Range("C214").FormatConditions.Add(1, 5, "15") | Out-Null
$ws.Range("B215").FormatConditions.Add(1, 5, "10") | Out-Null
$ws.Range("C215").FormatConditions.Add(1, 5, "15") | Out-Null
